$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for columns C (runs), D (balls), E (fours), F (sixes).
# Row 4 is intentionally omitted - its values are unchanged by the edit.
$data = @{
    2  = @(20, 12, 3, 1)
    3  = @(40, 25, 5, 1)
    5  = @(53, 43, 3, 3)
    6  = @(46, 37, 5, 2)
    7  = @(14, 15, 1, 0)
    8  = @(25, 13, 2, 2)
    9  = @(53, 36, 4, 3)
    10 = @(1, 3, 0, 0)
    11 = @(6, 4, 0, 1)
    12 = @(0, 5, 0, 0)
    13 = @(33, 20, 5, 0)
    14 = @(78, 44, 9, 3)
    15 = @(23, 15, 3, 1)
    16 = @(18, 19, 0, 1)
    17 = @(67, 39, 4, 4)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $targetRange = $ws.Range("C$row" + ":F$row")
    $targetRange.NumberFormat = "@"
    $ws.Cells.Item($row, 3).Value = [string]$vals[0]
    $ws.Cells.Item($row, 4).Value = [string]$vals[1]
    $ws.Cells.Item($row, 5).Value = [string]$vals[2]
    $ws.Cells.Item($row, 6).Value = [string]$vals[3]
}
